# Generate Report for Handoff
# Updates status text ("Handed back: in sync with en-US" -> "Ready for handoff")
# and the related handoff timestamps on all three sheets, plus the resulting
# narrower display width of the status/date columns.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-11-09 06:40:28"
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333336
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333336

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-11-09 06:40:15"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333336

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-11-09 06:40:28"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333336
